$d = $word.ActiveDocument

# --- 1. Insert two new completed checklist items right after the
#        "[x] Push branch feature/rider-backend-wiring ..." line (and before
#        the blank paragraph that precedes "Artifacts created:").
$anchorIndex = -1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Push branch feature/rider-backend-wiring from Windows*") {
        $anchorIndex = $p.Index
        break
    }
}

if ($anchorIndex -ge 1) {
    $anchor = $d.Paragraphs.Item($anchorIndex)
    $anchor.Range.InsertParagraphAfter()
    $d.Paragraphs.Item($anchorIndex + 1).Range.Text = "[x] Rider profile updates call backend (/users/me)"

    $d.Paragraphs.Item($anchorIndex + 1).Range.InsertParagraphAfter()
    $d.Paragraphs.Item($anchorIndex + 2).Range.Text = "[x] App config supports API base URL + mock toggle via dart-define"
}

# --- 2. Mark the "Choose vehicle type" checklist item as completed.
$d.Content.Find.Execute(
    "[~] (In Progress) Choose vehicle type (car/okada/pragya/aboboyaa)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[x] (Completed) Choose vehicle type (car/okada/pragya/aboboyaa)", 2)
